# Update the cryptos list with newly scraped prices / volume percentages.
# Rows 43 and 44 (Algorand / Aptos) also swap places in the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # The "price" column contains values that look numeric (e.g. "320.32",
    # "1.004"). A plain .Value assignment lets Excel auto-convert them to
    # real numbers (losing trailing zeros / formatting), so force the cell
    # to Text, assign the literal string, then restore the default style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 43 used to be Algorand, row 44 used to be Aptos.
# They swap order: row 43 becomes Aptos, row 44 becomes Algorand.
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D43") "10.28"
$ws.Range("E43").Value = "  -7.33%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D44") "0.1847"
$ws.Range("E44").Value = "  -3.86%  "

# Price (column D) and Volume(1h) percentage (column E) updates for every
# other row. Row 25 only has its volume percentage changed (price unchanged),
# and row 35 is untouched.
$updates = @(
    @{ Row = 2;  D = "27.445.16";   E = "  -3.45%  " },
    @{ Row = 3;  D = "1.851.76";    E = "  -5.02%  " },
    @{ Row = 4;  D = "1.004";       E = "  -0.95%  " },
    @{ Row = 5;  D = "320.32";      E = "  -0.26%  " },
    @{ Row = 6;  D = "1.005";       E = "  -0.63%  " },
    @{ Row = 7;  D = "0.4473";      E = "  -6.04%  " },
    @{ Row = 8;  D = "0.3846";      E = "  -5.04%  " },
    @{ Row = 9;  D = "47.72";       E = "  -10.66%  " },
    @{ Row = 10; D = "0.07828";     E = "  -7.65%  " },
    @{ Row = 11; D = "1.018";       E = "  -3.68%  " },
    @{ Row = 12; D = "21.34";       E = "  -3.33%  " },
    @{ Row = 13; D = "1.871.91";    E = "  -4.29%  " },
    @{ Row = 14; D = "5.858";       E = "  -5.09%  " },
    @{ Row = 15; D = "7.120";       E = "  -6.22%  " },
    @{ Row = 16; D = "1.009";       E = "  -0.48%  " },
    @{ Row = 17; D = "86.06";       E = "  -3.41%  " },
    @{ Row = 18; D = "0.00001025"; E = "  -4.61%  " },
    @{ Row = 19; D = "0.06510";     E = "  -1.33%  " },
    @{ Row = 20; D = "17.07";       E = "  -8.66%  " },
    @{ Row = 21; D = "1.004";       E = "  -0.85%  " },
    @{ Row = 22; D = "5.493";       E = "  -5.39%  " },
    @{ Row = 23; D = "27.429.04";   E = "  -3.59%  " },
    @{ Row = 24; D = "10.80";       E = "  -6.09%  " },
    @{ Row = 25; D = $null;         E = "  -0.31%  " },
    @{ Row = 26; D = "2.117.77";    E = "  -3.09%  " },
    @{ Row = 27; D = "150.30";      E = "  -2.66%  " },
    @{ Row = 28; D = "19.32";       E = "  -4.18%  " },
    @{ Row = 29; D = "5.505";       E = "  -7.55%  " },
    @{ Row = 30; D = "2.025";       E = "  -5.89%  " },
    @{ Row = 31; D = "120.19";      E = "  -2.85%  " },
    @{ Row = 32; D = "0.09376";     E = "  -1.94%  " },
    @{ Row = 33; D = "1.485";       E = "  +2.92%  " },
    @{ Row = 34; D = "0.9263";      E = "  -6.31%  " },
    @{ Row = 36; D = "5.246";       E = "  -6.28%  " },
    @{ Row = 37; D = "0.02221";     E = "  -4.94%  " },
    @{ Row = 38; D = "1.220";       E = "  -2.88%  " },
    @{ Row = 39; D = "0.05947";     E = "  -4.47%  " },
    @{ Row = 40; D = "8.318";       E = "  -5.32%  " },
    @{ Row = 41; D = "1.005";       E = "  -0.54%  " },
    @{ Row = 42; D = "0.5871";      E = "  -5.66%  " },
    @{ Row = 45; D = "1.283";       E = "  -3.66%  " },
    @{ Row = 46; D = "0.5624";      E = "  -5.69%  " },
    @{ Row = 47; D = "12.12";       E = "  -6.70%  " },
    @{ Row = 48; D = "1.918";       E = "  -6.71%  " },
    @{ Row = 49; D = "3.345";       E = "  -1.32%  " },
    @{ Row = 50; D = "0.06852";     E = "  +0.55%  " },
    @{ Row = 51; D = "1.047";       E = "  +3.39%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($r, 4) $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
